# Updates cryptos list figures (Price / Volume(1h)) to the latest scrape.
# Values are plain text (e.g. "44.872.75", "0.999", "  +1.25%  ") in the source
# sheet, so each cell is (re)formatted as Text before the write and reset back
# to the default "Normal" style afterwards -- this stops Excel's automatic
# "numbers typed as text get converted to numbers" behavior from mangling
# values like "7.20" -> 7.2 or "0.999" -> 0.999 (number) while leaving the
# cell style/formatting exactly as it was before the edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "45.089.83"
Set-TextValue "E2" "  +1.25%  "
Set-TextValue "D3" "2.265.39"
Set-TextValue "E3" "  +1.13%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.80%  "
Set-TextValue "D5" "301.34"
Set-TextValue "E5" "  -1.69%  "
Set-TextValue "D6" "94.22"
Set-TextValue "E6" "  -1.00%  "
Set-TextValue "D7" "0.566"
Set-TextValue "E7" "  -0.71%  "
Set-TextValue "E8" "  -0.60%  "
Set-TextValue "D9" "0.512"
Set-TextValue "E9" "  -1.48%  "
Set-TextValue "D10" "34.32"
Set-TextValue "E10" "  -1.76%  "
Set-TextValue "D11" "0.0788"
Set-TextValue "E11" "  -2.09%  "
Set-TextValue "D12" "7.20"
Set-TextValue "E12" "  -0.44%  "
Set-TextValue "D13" "0.104"
Set-TextValue "E13" "  -0.20%  "
Set-TextValue "D14" "2.611.78"
Set-TextValue "E14" "  +1.18%  "
Set-TextValue "D15" "2.266.15"
Set-TextValue "E15" "  -0.20%  "
Set-TextValue "D16" "13.72"
Set-TextValue "E16" "  +0.98%  "
Set-TextValue "D17" "0.799"
Set-TextValue "E17" "  -4.34%  "
Set-TextValue "D18" "44.965.33"
Set-TextValue "E18" "  +1.55%  "
Set-TextValue "D19" "12.93"
Set-TextValue "E19" "  +9.01%  "
Set-TextValue "D20" "0.0₃0922"
Set-TextValue "E20" "  -3.23%  "
Set-TextValue "D21" "6.09"
Set-TextValue "E21" "  -3.43%  "
Set-TextValue "D22" "65.29"
Set-TextValue "E22" "  -0.25%  "
Set-TextValue "D23" "240.20"
Set-TextValue "E23" "  +1.27%  "
Set-TextValue "D24" "2.89"
Set-TextValue "E24" "  -2.33%  "
Set-TextValue "D25" "0.999"
Set-TextValue "E25" "  -0.50%  "
Set-TextValue "D26" "1.91"
Set-TextValue "E26" "  -3.70%  "
Set-TextValue "D27" "40.02"
Set-TextValue "E27" "  +6.81%  "
Set-TextValue "D28" "2.29"
Set-TextValue "E28" "  +0.37%  "
Set-TextValue "D29" "9.56"
Set-TextValue "E29" "  -2.35%  "
Set-TextValue "D30" "19.60"
Set-TextValue "E30" "  -1.68%  "
Set-TextValue "D31" "153.23"
Set-TextValue "E31" "  +0.53%  "
Set-TextValue "D32" "5.59"
Set-TextValue "E32" "  -5.96%  "
Set-TextValue "D33" "0.0793"
Set-TextValue "E33" "  -0.36%  "
Set-TextValue "D34" "2.57"
Set-TextValue "E34" "  -2.50%  "
Set-TextValue "E35" "  -1.06%  "
# Row 36: coin swapped into this slot
Set-TextValue "B36" "LidoDAOToken"
Set-TextValue "C36" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D36" "2.91"
Set-TextValue "E36" "  -4.36%  "
# Row 37: coin swapped into this slot
Set-TextValue "B37" "Kaspa"
Set-TextValue "C37" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D37" "0.106"
Set-TextValue "E37" "  -3.29%  "
Set-TextValue "D38" "1.74"
Set-TextValue "E38" "  -5.12%  "
Set-TextValue "D39" "0.0305"
Set-TextValue "E39" "  +1.37%  "
# Row 40: coin swapped into this slot
Set-TextValue "B40" "NEARProtocol"
Set-TextValue "C40" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D40" "3.28"
Set-TextValue "E40" "  -3.00%  "
# Row 41: coin swapped into this slot
Set-TextValue "B41" "RenderToken"
Set-TextValue "C41" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D41" "3.74"
Set-TextValue "E41" "  -0.72%  "
Set-TextValue "D42" "13.91"
Set-TextValue "E42" "  -6.78%  "
Set-TextValue "E43" "  -0.87%  "
Set-TextValue "D44" "1.785.94"
Set-TextValue "E44" "  -1.44%  "
Set-TextValue "D45" "1.85"
Set-TextValue "E45" "  +8.73%  "
Set-TextValue "E46" "  +0.45%  "
Set-TextValue "D47" "70.59"
Set-TextValue "E47" "  -0.05%  "
Set-TextValue "D48" "75.74"
Set-TextValue "E48" "  -3.77%  "
Set-TextValue "D49" "96.59"
Set-TextValue "E49" "  -2.29%  "
Set-TextValue "D50" "4.69"
Set-TextValue "E50" "  -4.28%  "
Set-TextValue "D51" "7.85"
Set-TextValue "E51" "  -2.14%  "
